$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, pushing the existing rows 8-14 down to 9-15
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the weekly price record
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44915
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100114007
$ws.Range("G8").Value = "Jengibre"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 18000
$ws.Range("N8").Value = "`$/caja 13 kilos"
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 1385
$ws.Range("Q8").Value = 13
$ws.Range("R8").Value = "Hortaliza"
